$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Productdata": update StartingInventories (column C) and
# SetupCosts (column E) for rows 2..23.
# ---------------------------------------------------------------------------
$wsProd = $wb.Worksheets.Item("Productdata")

$prodC = @{
    2 = 0;     3 = 9001;  4 = 2251;  5 = 0;     6 = 4501;  7 = 1351;
    8 = 451;   9 = 6301;  10 = 4501; 11 = 13501; 12 = 3151; 13 = 0;
    14 = 0;    15 = 0;    16 = 0;    17 = 0;    18 = 0;    19 = 0;
    20 = 1;    21 = 1;    22 = 1;    23 = 1
}

$prodE = @{
    2 = 192.456;             3 = 69.92999999999999;   4 = 17.50333333333333;
    5 = 24.37516666666667;   6 = 34.16333333333333;   7 = 10.272;
    8 = 3.409333333333333;   9 = 50.49566666666666;   10 = 35.99999999999999;
    11 = 107.445;            12 = 25.235;              13 = 184.338;
    14 = 68.41666666666666;  15 = 16.87416666666667;  16 = 23.37766666666666;
    17 = 33.49499999999999;  18 = 10.071;              19 = 3.343333333333333;
    20 = 43.83333333333334;  21 = 46.53333333333333;  22 = 57.59999999999999;
    23 = 177.2333333333333
}

foreach ($row in 2..23) {
    $wsProd.Range("C$row").Value = $prodC[$row]
    $wsProd.Range("E$row").Value = $prodE[$row]
}

# ---------------------------------------------------------------------------
# Sheet "Capacity": update column B for rows 2..23.
# ---------------------------------------------------------------------------
$wsCap = $wb.Worksheets.Item("Capacity")

$capB = @{
    2 = 48600;  3 = 18000;  4 = 1500;   5 = 8400;   6 = 12000;  7 = 900;
    8 = 1500;   9 = 21000;  10 = 15000; 11 = 36000; 12 = 4200;  13 = 48600;
    14 = 6000;  15 = 7500;  16 = 4200;  17 = 3000;  18 = 1800;  19 = 600;
    20 = 120000; 21 = 60000; 22 = 90000; 23 = 90000
}

foreach ($row in 2..23) {
    $wsCap.Range("B$row").Value = $capB[$row]
}

# ---------------------------------------------------------------------------
# Sheet "ProcessingTime": update the non-zero diagonal cells.
# ---------------------------------------------------------------------------
$wsProc = $wb.Worksheets.Item("ProcessingTime")

$procCells = @{
    "B2" = 3;  "D4" = 1;  "E5" = 4;  "F6" = 4;  "G7" = 1;  "H8" = 5;
    "I9" = 5;  "J10" = 5; "L12" = 2; "M13" = 3; "N14" = 1; "O15" = 5;
    "P16" = 2; "R18" = 2; "S19" = 2; "T20" = 4; "U21" = 2; "V22" = 3;
    "W23" = 3
}

foreach ($addr in $procCells.Keys) {
    $wsProc.Range($addr).Value = $procCells[$addr]
}
